$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$donorStyle1 = $ws.Range("A1")
$donorStyle2 = $ws.Range("D2")
$donorStyle3 = $ws.Range("N2")

# Match the existing column width (raw OOXML width="12") used by every other
# data column so the new CQ column (95) renders identically.
$ws.Range("CQ1").ColumnWidth = 11.17

# Row 1: style 1, value '2024/12/12'
$ws.Range("CQ1").NumberFormat = "@"
$ws.Range("CQ1").Value = "2024/12/12"
$donorStyle1.Copy()
$ws.Range("CQ1").PasteSpecial(-4122)

# Row 2: style 1, value 186.3
$ws.Range("CQ2").Value = 186.3
$donorStyle1.Copy()
$ws.Range("CQ2").PasteSpecial(-4122)
$ws.Range("CQ2").Value = 186.3

# Row 3: style 1, value 192.4
$ws.Range("CQ3").Value = 192.4
$donorStyle1.Copy()
$ws.Range("CQ3").PasteSpecial(-4122)
$ws.Range("CQ3").Value = 192.4

# Row 4: style 1, value 188.2
$ws.Range("CQ4").Value = 188.2
$donorStyle1.Copy()
$ws.Range("CQ4").PasteSpecial(-4122)
$ws.Range("CQ4").Value = 188.2

# Row 5: style 3, value 136.9
$ws.Range("CQ5").Value = 136.9
$donorStyle3.Copy()
$ws.Range("CQ5").PasteSpecial(-4122)
$ws.Range("CQ5").Value = 136.9

# Row 6: style 3, value 128.7
$ws.Range("CQ6").Value = 128.7
$donorStyle3.Copy()
$ws.Range("CQ6").PasteSpecial(-4122)
$ws.Range("CQ6").Value = 128.7

# Row 7: style 2, value 114
$ws.Range("CQ7").Value = 114
$donorStyle2.Copy()
$ws.Range("CQ7").PasteSpecial(-4122)
$ws.Range("CQ7").Value = 114

# Row 8: style 3, value 127.6
$ws.Range("CQ8").Value = 127.6
$donorStyle3.Copy()
$ws.Range("CQ8").PasteSpecial(-4122)
$ws.Range("CQ8").Value = 127.6

# Row 9: style 2, value 124
$ws.Range("CQ9").Value = 124
$donorStyle2.Copy()
$ws.Range("CQ9").PasteSpecial(-4122)
$ws.Range("CQ9").Value = 124

# Row 10: style 1, value 226.2
$ws.Range("CQ10").Value = 226.2
$donorStyle1.Copy()
$ws.Range("CQ10").PasteSpecial(-4122)
$ws.Range("CQ10").Value = 226.2

# Row 11: style 1, value 176.1
$ws.Range("CQ11").Value = 176.1
$donorStyle1.Copy()
$ws.Range("CQ11").PasteSpecial(-4122)
$ws.Range("CQ11").Value = 176.1

# Row 12: style 1, value 143.6
$ws.Range("CQ12").Value = 143.6
$donorStyle1.Copy()
$ws.Range("CQ12").PasteSpecial(-4122)
$ws.Range("CQ12").Value = 143.6

# Row 13: style 2, value 124.4
$ws.Range("CQ13").Value = 124.4
$donorStyle2.Copy()
$ws.Range("CQ13").PasteSpecial(-4122)
$ws.Range("CQ13").Value = 124.4

# Row 14: style 3, value 136.7
$ws.Range("CQ14").Value = 136.7
$donorStyle3.Copy()
$ws.Range("CQ14").PasteSpecial(-4122)
$ws.Range("CQ14").Value = 136.7

# Row 15: style 1, value 157.2
$ws.Range("CQ15").Value = 157.2
$donorStyle1.Copy()
$ws.Range("CQ15").PasteSpecial(-4122)
$ws.Range("CQ15").Value = 157.2

# Row 16: style 1, value 158.5
$ws.Range("CQ16").Value = 158.5
$donorStyle1.Copy()
$ws.Range("CQ16").PasteSpecial(-4122)
$ws.Range("CQ16").Value = 158.5

# Row 17: style 1, value 157.2
$ws.Range("CQ17").Value = 157.2
$donorStyle1.Copy()
$ws.Range("CQ17").PasteSpecial(-4122)
$ws.Range("CQ17").Value = 157.2

# Row 18: style 3, value 137.7
$ws.Range("CQ18").Value = 137.7
$donorStyle3.Copy()
$ws.Range("CQ18").PasteSpecial(-4122)
$ws.Range("CQ18").Value = 137.7

# Row 19: style 3, value 139.8
$ws.Range("CQ19").Value = 139.8
$donorStyle3.Copy()
$ws.Range("CQ19").PasteSpecial(-4122)
$ws.Range("CQ19").Value = 139.8

# Row 20: style 1, value 150.5
$ws.Range("CQ20").Value = 150.5
$donorStyle1.Copy()
$ws.Range("CQ20").PasteSpecial(-4122)
$ws.Range("CQ20").Value = 150.5

# Row 21: style 3, value 136.9
$ws.Range("CQ21").Value = 136.9
$donorStyle3.Copy()
$ws.Range("CQ21").PasteSpecial(-4122)
$ws.Range("CQ21").Value = 136.9

# Row 22: style 1, value 168.5
$ws.Range("CQ22").Value = 168.5
$donorStyle1.Copy()
$ws.Range("CQ22").PasteSpecial(-4122)
$ws.Range("CQ22").Value = 168.5

# Row 23: style 1, value 148.1
$ws.Range("CQ23").Value = 148.1
$donorStyle1.Copy()
$ws.Range("CQ23").PasteSpecial(-4122)
$ws.Range("CQ23").Value = 148.1

# Row 24: style 1, value 152.8
$ws.Range("CQ24").Value = 152.8
$donorStyle1.Copy()
$ws.Range("CQ24").PasteSpecial(-4122)
$ws.Range("CQ24").Value = 152.8

# Row 25: style 2, value 111.2
$ws.Range("CQ25").Value = 111.2
$donorStyle2.Copy()
$ws.Range("CQ25").PasteSpecial(-4122)
$ws.Range("CQ25").Value = 111.2

# Row 26: style 1, value 186.9
$ws.Range("CQ26").Value = 186.9
$donorStyle1.Copy()
$ws.Range("CQ26").PasteSpecial(-4122)
$ws.Range("CQ26").Value = 186.9

# Row 27: style 1, value 160.3
$ws.Range("CQ27").Value = 160.3
$donorStyle1.Copy()
$ws.Range("CQ27").PasteSpecial(-4122)
$ws.Range("CQ27").Value = 160.3

# Row 28: style 1, value 178.9
$ws.Range("CQ28").Value = 178.9
$donorStyle1.Copy()
$ws.Range("CQ28").PasteSpecial(-4122)
$ws.Range("CQ28").Value = 178.9

# Row 29: style 1, value 162.3
$ws.Range("CQ29").Value = 162.3
$donorStyle1.Copy()
$ws.Range("CQ29").PasteSpecial(-4122)
$ws.Range("CQ29").Value = 162.3

# Row 30: style 3, value 132.4
$ws.Range("CQ30").Value = 132.4
$donorStyle3.Copy()
$ws.Range("CQ30").PasteSpecial(-4122)
$ws.Range("CQ30").Value = 132.4

# Row 31: style 1, value 183.6
$ws.Range("CQ31").Value = 183.6
$donorStyle1.Copy()
$ws.Range("CQ31").PasteSpecial(-4122)
$ws.Range("CQ31").Value = 183.6

# Row 32: style 1, value 200.5
$ws.Range("CQ32").Value = 200.5
$donorStyle1.Copy()
$ws.Range("CQ32").PasteSpecial(-4122)
$ws.Range("CQ32").Value = 200.5

# Row 33: style 3, value 128.8
$ws.Range("CQ33").Value = 128.8
$donorStyle3.Copy()
$ws.Range("CQ33").PasteSpecial(-4122)
$ws.Range("CQ33").Value = 128.8

# Row 34: style 1, value 219.8
$ws.Range("CQ34").Value = 219.8
$donorStyle1.Copy()
$ws.Range("CQ34").PasteSpecial(-4122)
$ws.Range("CQ34").Value = 219.8

# Row 35: style 3, value 132
$ws.Range("CQ35").Value = 132
$donorStyle3.Copy()
$ws.Range("CQ35").PasteSpecial(-4122)
$ws.Range("CQ35").Value = 132

# Row 36: style 1, value 173.8
$ws.Range("CQ36").Value = 173.8
$donorStyle1.Copy()
$ws.Range("CQ36").PasteSpecial(-4122)
$ws.Range("CQ36").Value = 173.8

# Row 37: style 1, value 165.2
$ws.Range("CQ37").Value = 165.2
$donorStyle1.Copy()
$ws.Range("CQ37").PasteSpecial(-4122)
$ws.Range("CQ37").Value = 165.2

# Row 38: style 3, value 135.2
$ws.Range("CQ38").Value = 135.2
$donorStyle3.Copy()
$ws.Range("CQ38").PasteSpecial(-4122)
$ws.Range("CQ38").Value = 135.2

# Row 39: style 1, value 214.4
$ws.Range("CQ39").Value = 214.4
$donorStyle1.Copy()
$ws.Range("CQ39").PasteSpecial(-4122)
$ws.Range("CQ39").Value = 214.4

# Row 40: style 1, value 140.9
$ws.Range("CQ40").Value = 140.9
$donorStyle1.Copy()
$ws.Range("CQ40").PasteSpecial(-4122)
$ws.Range("CQ40").Value = 140.9

# Row 41: style 3, value 130.1
$ws.Range("CQ41").Value = 130.1
$donorStyle3.Copy()
$ws.Range("CQ41").PasteSpecial(-4122)
$ws.Range("CQ41").Value = 130.1

# Row 42: style 1, value 162.1
$ws.Range("CQ42").Value = 162.1
$donorStyle1.Copy()
$ws.Range("CQ42").PasteSpecial(-4122)
$ws.Range("CQ42").Value = 162.1

# Row 43: style 1, value 141.4
$ws.Range("CQ43").Value = 141.4
$donorStyle1.Copy()
$ws.Range("CQ43").PasteSpecial(-4122)
$ws.Range("CQ43").Value = 141.4

# Row 44: style 1, value 147.1
$ws.Range("CQ44").Value = 147.1
$donorStyle1.Copy()
$ws.Range("CQ44").PasteSpecial(-4122)
$ws.Range("CQ44").Value = 147.1

# Row 45: style 3, value 131.8
$ws.Range("CQ45").Value = 131.8
$donorStyle3.Copy()
$ws.Range("CQ45").PasteSpecial(-4122)
$ws.Range("CQ45").Value = 131.8

# Row 46: style 1, value 146.8
$ws.Range("CQ46").Value = 146.8
$donorStyle1.Copy()
$ws.Range("CQ46").PasteSpecial(-4122)
$ws.Range("CQ46").Value = 146.8

# Row 47: style 3, value 131.6
$ws.Range("CQ47").Value = 131.6
$donorStyle3.Copy()
$ws.Range("CQ47").PasteSpecial(-4122)
$ws.Range("CQ47").Value = 131.6

# Row 48: style 1, value 162.3
$ws.Range("CQ48").Value = 162.3
$donorStyle1.Copy()
$ws.Range("CQ48").PasteSpecial(-4122)
$ws.Range("CQ48").Value = 162.3

# Row 49: style 1, value 248.2
$ws.Range("CQ49").Value = 248.2
$donorStyle1.Copy()
$ws.Range("CQ49").PasteSpecial(-4122)
$ws.Range("CQ49").Value = 248.2

# Row 50: style 1, value 183.4
$ws.Range("CQ50").Value = 183.4
$donorStyle1.Copy()
$ws.Range("CQ50").PasteSpecial(-4122)
$ws.Range("CQ50").Value = 183.4

# Row 51: style 1, value 221.4
$ws.Range("CQ51").Value = 221.4
$donorStyle1.Copy()
$ws.Range("CQ51").PasteSpecial(-4122)
$ws.Range("CQ51").Value = 221.4

# Row 52: style 3, value 134.9
$ws.Range("CQ52").Value = 134.9
$donorStyle3.Copy()
$ws.Range("CQ52").PasteSpecial(-4122)
$ws.Range("CQ52").Value = 134.9

# Row 53: style 1, value 142.8
$ws.Range("CQ53").Value = 142.8
$donorStyle1.Copy()
$ws.Range("CQ53").PasteSpecial(-4122)
$ws.Range("CQ53").Value = 142.8

$excel.CutCopyMode = $false
Write-Output "Added 2024/12/12 column (CQ) with 52 data rows"